$wb = $excel.ActiveWorkbook

# --- Section_A ---
$ws = $wb.Worksheets.Item("Section_A")
$ws.Range("B2").Value = "CS262"
$ws.Range("C2").Value = "HS261 (Elective)"
$ws.Range("D2").Value = "MA262"
$ws.Range("E2").Value = "HS261 (Elective)"
$ws.Range("F2").Value = "CS263"
$ws.Range("B3").Value = "CS263"
$ws.Range("C3").Value = "MA261"
$ws.Range("D3").Value = "CS262"
$ws.Range("E3").Value = "CS264"
$ws.Range("C5").Value = "CS263"
$ws.Range("D5").Value = "CS251 (Elective)"
$ws.Range("E5").Value = "CS261"
$ws.Range("F5").Value = "MA262"
$ws.Range("B6").Value = "HS261 (Tutorial)"
$ws.Range("C6").Value = "Free"
$ws.Range("D6").Value = "Free"
$ws.Range("F6").Value = "CS264 (Tutorial)"
$ws.Range("B7").Value = "CS261"
$ws.Range("C7").Value = "CS261"
$ws.Range("D7").Value = "CS264"
$ws.Range("E7").Value = "MA261"
$ws.Range("D8").Value = "Free"
$ws.Range("E8").Value = "Free"
$ws.Range("F8").Value = "CS251 (Tutorial)"

# --- Section_B ---
$ws = $wb.Worksheets.Item("Section_B")
$ws.Range("B2").Value = "MA261"
$ws.Range("C2").Value = "HS261 (Elective)"
$ws.Range("D2").Value = "CS263"
$ws.Range("E2").Value = "HS261 (Elective)"
$ws.Range("F2").Value = "CS261"
$ws.Range("B3").Value = "CS263"
$ws.Range("C3").Value = "CS264"
$ws.Range("D3").Value = "CS261"
$ws.Range("E3").Value = "CS263"
$ws.Range("F3").Value = "CS264"
$ws.Range("C5").Value = "MA262"
$ws.Range("D5").Value = "CS251 (Elective)"
$ws.Range("E5").Value = "CS264"
$ws.Range("F5").Value = "Free"
$ws.Range("B6").Value = "HS261 (Tutorial)"
$ws.Range("C6").Value = "Free"
$ws.Range("D6").Value = "Free"
$ws.Range("F6").Value = "Free"
$ws.Range("B7").Value = "CS261"
$ws.Range("C7").Value = "CS262"
$ws.Range("D7").Value = "CS262"
$ws.Range("E7").Value = "MA262"
$ws.Range("F7").Value = "MA261"
$ws.Range("B8").Value = "CS264 (Tutorial)"
$ws.Range("C8").Value = "Free"
$ws.Range("D8").Value = "Free"
$ws.Range("E8").Value = "Free"
$ws.Range("F8").Value = "CS251 (Tutorial)"

# --- Elective_Coordination ---
$ws = $wb.Worksheets.Item("Elective_Coordination")
$ws.Range("C2").Value = "Tue"
$ws.Range("C3").Value = "Thu"
$ws.Range("C4").Value = "Mon"
$ws.Range("D4").Value = "14:30-15:30"
$ws.Range("C11").Value = "Mon"
$ws.Range("C12").Value = "Wed"
$ws.Range("C13").Value = "Fri"
$ws.Range("D13").Value = "17:00-18:00"
